# The "falling" gyroscope sample was re-captured: two fresh x/y/z samples
# were prepended to the series (rows 2-3) and every remaining sample shifted
# down by two rows. Because the series keeps a fixed-length rolling window,
# the oldest sample (the former last row, r=22) drops off the sheet, so the
# sheet shrinks from A1:C22 to A1:C21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full x,y,z gyroscope series as it should look after the re-capture,
# written starting at row 2 (row 1 holds the "x","y","z" headers).
$data = @(
    @(-0.06746154287180267,  0.03214145534087062, -0.9956535080144557),
    @( 0.2888078393284657,   0.04184071745636832, -1.758943790613218),
    @( 1.293112953041876,   -0.7586058728223629,  -2.201892187429029),
    @( 1.076085335986568,   -0.3218408019043786,  -1.197818033223937),
    @( 1.466690946456983,   -0.5813377291657031,  -0.9997129280899739),
    @( 4.411398765652687,   -4.482220122980524,   -3.045610698849662),
    @( 0.08928575072178435, -1.550635254660315,    0.06088052863297655),
    @(-0.9927981232487877,  -1.707755399304766,    1.41453302322432),
    @( 2.431215341007984,    1.645354886387669,    0.273621466270713),
    @( 0.2498580986677279,   0.9151578448539578,   0.4712292718332862),
    @( 0.04947298288691901,  0.4319243583568316,   0.3634221997371941),
    @(-0.5032464119410757,  -0.05431018042009644,  0.3147802685582363),
    @(-0.8943670974221336,  -1.390350583680838,    2.061716448429016),
    @( 0.08820253887841945, -0.5515473787347926,   0.8748975694179519),
    @( 0.2351689076181083,  -0.03523484912029523,  0.06333464960199484),
    @(-0.005455169480207168,-0.04037748242533482,  0.09177540264330601),
    @(-0.03705323929356984, -0.02731135879578294,  0.08058094122815343),
    @(-0.03103692665003059,  0.001527163070128304, 0.01256180448500908),
    @(-0.01491292483758091, -0.00055048926625147, -0.0008630249216113467),
    @( 0.02735042602343611, -0.0430127267920694,  -0.02982940334220257)
)

$row = 2
foreach ($sample in $data) {
    $ws.Range("A$row").Value = $sample[0]
    $ws.Range("B$row").Value = $sample[1]
    $ws.Range("C$row").Value = $sample[2]
    $row++
}

# Drop the now-stale trailing sample so the sheet ends at row 21 (A1:C21).
$ws.Rows(22).Delete()
